# "few changes in deals page"
# - rename "Sheet1" to "deal" and populate it with deal data
# - re-point the active tab/selection to the new deal sheet
# - tweak the "contact" sheet's remembered selection/zoom

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("contact")
$deal = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 -> deal
$deal.Name = "deal"

# Populate header row (write in this column order so new shared-string
# indices line up: Pipeline, Amount, DealName, DealStage)
$deal.Range("B1").Value = "Pipeline"
$deal.Range("D1").Value = "Amount"
$deal.Range("A1").Value = "DealName"
$deal.Range("C1").Value = "DealStage"

# Row 2
$deal.Range("A2").Value = "deal1"
$deal.Range("B2").Value = "pipe1"
$deal.Range("C2").Value = "Appointment Scheduled"
$deal.Range("D2").Value = 15000

# Row 3
$deal.Range("A3").Value = "deal2"
$deal.Range("B3").Value = "pipe2"
$deal.Range("C3").Value = "Qualified To Buy"
$deal.Range("D3").Value = 21000

# Header row styling: yellow fill, regular (non-bold) font
$deal.Range("A1:D1").Interior.Color = 65535

# Widen the DealStage column so the longest label fits
$deal.Columns.Item(3).ColumnWidth = 19.3

# Move the "contact" sheet's remembered selection/zoom
$contact.Range("B13").Select()
$excel.ActiveWindow.Zoom = 66

# Make "deal" the active/visible sheet with its own remembered selection
$deal.Activate()
$deal.Range("C11").Select()
